$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hindcast")
$ws.Range("B2").Value2 = 184408164.9387597
$ws.Range("C2").Value2 = 92259838.91505134
$ws.Range("B3").Value2 = 105776518.8397668
$ws.Range("C3").Value2 = 52115801.83860195
$ws.Range("B4").Value2 = 65319688.74451113
$ws.Range("C4").Value2 = 29915697.52687296
$ws.Range("B5").Value2 = 52618397.08947015
$ws.Range("C5").Value2 = 23218493.31353269
$ws.Range("B6").Value2 = 52393409.30959573
$ws.Range("C6").Value2 = 22593164.83046902
$ws.Range("B7").Value2 = 43913787.16428258
$ws.Range("C7").Value2 = 18806449.29957822
$ws.Range("B8").Value2 = 41677204.76110819
$ws.Range("C8").Value2 = 17923062.95776032
$ws.Range("B9").Value2 = 48923285.39682809
$ws.Range("C9").Value2 = 22136835.62842614
$ws.Range("B10").Value2 = 75137807.59089747
$ws.Range("C10").Value2 = 35863840.94166603
$ws.Range("B11").Value2 = 106801027.5808776
$ws.Range("C11").Value2 = 50456388.3182405
$ws.Range("B12").Value2 = 112104697.4033566
$ws.Range("C12").Value2 = 53423244.57584248
$ws.Range("B13").Value2 = 126389102.6321382
$ws.Range("C13").Value2 = 62891365.8917788
$ws.Range("B14").Value2 = 142321291.3542317
$ws.Range("C14").Value2 = 74814480.68880057

$ws = $wb.Worksheets.Item("condensed")
$ws.Range("B2").Value2 = 193279770.1022453
$ws.Range("C2").Value2 = 89528517.06589615
$ws.Range("B3").Value2 = 110915108.5467549
$ws.Range("C3").Value2 = 52292898.22731744
$ws.Range("B4").Value2 = 67816337.12785231
$ws.Range("C4").Value2 = 28691080.63697247
$ws.Range("B5").Value2 = 54652968.22227819
$ws.Range("C5").Value2 = 22195065.09432386
$ws.Range("B6").Value2 = 49848928.55053227
$ws.Range("C6").Value2 = 19741346.80022854
$ws.Range("B7").Value2 = 34535974.95571143
$ws.Range("C7").Value2 = 13881193.43158263
$ws.Range("B8").Value2 = 36441826.52468979
$ws.Range("C8").Value2 = 14704081.66305494
$ws.Range("B9").Value2 = 48260907.26081733
$ws.Range("C9").Value2 = 20110258.02607381
$ws.Range("B10").Value2 = 75223092.23542337
$ws.Range("C10").Value2 = 33181543.61009645
$ws.Range("B11").Value2 = 108905480.1157428
$ws.Range("C11").Value2 = 47881842.16304622
$ws.Range("B12").Value2 = 115628944.8802835
$ws.Range("C12").Value2 = 51805955.58964063
$ws.Range("B13").Value2 = 132469216.9360885
$ws.Range("C13").Value2 = 63175849.86661611
$ws.Range("B14").Value2 = 148907182.5884926
$ws.Range("C14").Value2 = 76278705.40769833
